$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Cell text substitutions (using original row/col indices, before insert/delete) ---
$t.Cell(1,1).Range.Text = "27+18=45"
$t.Cell(1,2).Range.Text = "64-62=2"
$t.Cell(1,3).Range.Text = "71+22=93"
$t.Cell(1,4).Range.Text = "6+53=59"
$t.Cell(1,5).Range.Text = "80+9=89"
$t.Cell(2,1).Range.Text = "18+3=21"
$t.Cell(2,2).Range.Text = "17+7=24"
$t.Cell(2,3).Range.Text = "95-51=44"
$t.Cell(2,4).Range.Text = "13+32=45"
$t.Cell(2,5).Range.Text = "14+57=71"
$t.Cell(3,1).Range.Text = "2+84=86"
$t.Cell(3,2).Range.Text = "98-13=85"
$t.Cell(3,3).Range.Text = "48-33=15"
$t.Cell(3,4).Range.Text = "10+63=73"
$t.Cell(3,5).Range.Text = "9+67=76"
$t.Cell(4,1).Range.Text = "17+18=35"
$t.Cell(4,2).Range.Text = "69+2=71"
$t.Cell(4,3).Range.Text = "66-66=0"
$t.Cell(4,4).Range.Text = "77-3=74"
$t.Cell(4,5).Range.Text = "66-54=12"
$t.Cell(5,1).Range.Text = "61+5=66"
$t.Cell(5,2).Range.Text = "22+66=88"
$t.Cell(5,3).Range.Text = "3+67=70"
$t.Cell(5,4).Range.Text = "93-69=24"
$t.Cell(5,5).Range.Text = "18+31=49"
$t.Cell(6,1).Range.Text = "82-27=55"
$t.Cell(6,2).Range.Text = "15+70=85"
$t.Cell(6,3).Range.Text = "63+28=91"
$t.Cell(6,4).Range.Text = "57+12=69"
$t.Cell(6,5).Range.Text = "98-69=29"
$t.Cell(7,1).Range.Text = "72-5=67"
$t.Cell(7,2).Range.Text = "12+45=57"
$t.Cell(7,3).Range.Text = "17+63=80"
$t.Cell(7,4).Range.Text = "66+31=97"
$t.Cell(7,5).Range.Text = "26+33=59"
$t.Cell(8,1).Range.Text = "63-38=25"
$t.Cell(8,2).Range.Text = "60+26=86"
$t.Cell(8,3).Range.Text = "74-39=35"
$t.Cell(8,4).Range.Text = "61-25=36"
$t.Cell(8,5).Range.Text = "74-53=21"
$t.Cell(9,1).Range.Text = "19+49=68"
$t.Cell(9,2).Range.Text = "63+15=78"
$t.Cell(9,3).Range.Text = "33+51=84"
$t.Cell(9,4).Range.Text = "40+2=42"
$t.Cell(9,5).Range.Text = "9+65=74"
$t.Cell(10,2).Range.Text = "23-5=18"
$t.Cell(10,3).Range.Text = "71-38=33"
$t.Cell(10,4).Range.Text = "18+22=40"
$t.Cell(10,5).Range.Text = "71-51=20"
$t.Cell(11,1).Range.Text = "17+0=17"
$t.Cell(11,2).Range.Text = "9+54=63"
$t.Cell(11,3).Range.Text = "78-12=66"
$t.Cell(11,4).Range.Text = "54-12=42"
$t.Cell(11,5).Range.Text = "22-13=9"
$t.Cell(12,1).Range.Text = "26+5=31"
$t.Cell(12,2).Range.Text = "25+15=40"
$t.Cell(12,3).Range.Text = "54-5=49"
$t.Cell(12,4).Range.Text = "20-16=4"
$t.Cell(12,5).Range.Text = "13+43=56"
$t.Cell(13,1).Range.Text = "89-40=49"
$t.Cell(13,2).Range.Text = "24+7=31"
$t.Cell(13,3).Range.Text = "69-35=34"
$t.Cell(13,4).Range.Text = "65-19=46"
$t.Cell(13,5).Range.Text = "49+43=92"
$t.Cell(14,1).Range.Text = "7+71=78"
$t.Cell(14,2).Range.Text = "76-60=16"
$t.Cell(14,3).Range.Text = "16+31=47"
$t.Cell(14,4).Range.Text = "97-57=40"
$t.Cell(14,5).Range.Text = "85-1=84"
$t.Cell(15,1).Range.Text = "12+7=19"
$t.Cell(15,2).Range.Text = "47-24=23"
$t.Cell(15,3).Range.Text = "26-1=25"
$t.Cell(15,4).Range.Text = "84+3=87"
$t.Cell(15,5).Range.Text = "39+56=95"
$t.Cell(16,1).Range.Text = "10+79=89"
$t.Cell(16,2).Range.Text = "43+33=76"
$t.Cell(16,3).Range.Text = "34+59=93"
$t.Cell(16,4).Range.Text = "13-9=4"
$t.Cell(16,5).Range.Text = "37+33=70"
$t.Cell(17,1).Range.Text = "28+23=51"
$t.Cell(17,2).Range.Text = "90-63=27"
$t.Cell(17,3).Range.Text = "62+33=95"
$t.Cell(17,4).Range.Text = "34+32=66"
$t.Cell(17,5).Range.Text = "61-33=28"
$t.Cell(18,1).Range.Text = "64-47=17"
$t.Cell(18,2).Range.Text = "86-52=34"
$t.Cell(18,3).Range.Text = "13+53=66"
$t.Cell(18,4).Range.Text = "7+38=45"
$t.Cell(18,5).Range.Text = "93-72=21"

# --- Insert two new rows after original row 9 (before original row 10) ---
$beforeRow = $t.Rows.Item(10)
$newRowA = $t.Rows.Add($beforeRow)
$newRowA.Cells.Item(1).Range.Text = "82-81=1"
$newRowA.Cells.Item(2).Range.Text = "96-35=61"
$newRowA.Cells.Item(3).Range.Text = "62+37=99"
$newRowA.Cells.Item(4).Range.Text = "24+21=45"
$newRowA.Cells.Item(5).Range.Text = "60-47=13"
$beforeRow2 = $t.Rows.Item(11)
$newRowB = $t.Rows.Add($beforeRow2)
$newRowB.Cells.Item(1).Range.Text = "52-6=46"
$newRowB.Cells.Item(2).Range.Text = "78-48=30"
$newRowB.Cells.Item(3).Range.Text = "23-18=5"
$newRowB.Cells.Item(4).Range.Text = "87-13=74"
$newRowB.Cells.Item(5).Range.Text = "43-15=28"

# --- Delete the last two rows (originally rows 19 and 20, now shifted by +2 = 21,22) ---
$total = $t.Rows.Count
$t.Rows.Item($total).Delete()
$t.Rows.Item($total - 1).Delete()
